$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 17250
$ws.Range("I70").Value = 4500
$ws.Range("K70").Value = 13500
$ws.Range("M70").Value = -13230
$ws.Range("H73").Value = 17250
$ws.Range("I73").Value = 4500
$ws.Range("K73").Value = 13500
$ws.Range("M73").Value = -12564
$ws.Range("H112").Value = 1816.569
$ws.Range("J112").Value = 1844.7637
$ws.Range("L112").Value = 5534.2911
$ws.Range("N112").Value = -7750.2911
$ws.Range("H116").Value = 11083.272
$ws.Range("I116").Value = 18546.666
$ws.Range("J116").Value = 2127.2
$ws.Range("K116").Value = 18546.666
$ws.Range("L116").Value = 2127.2
$ws.Range("M116").Value = -15104.666
$ws.Range("N116").Value = -9011.200000000001
$ws.Range("H129").Value = 1237.5769
$ws.Range("I129").Value = 377
$ws.Range("J129").Value = 1272
$ws.Range("K129").Value = 1131
$ws.Range("L129").Value = 3816
$ws.Range("M129").Value = 3869
$ws.Range("N129").Value = -13816

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5448.1523
$ws.Range("I32").Value = 3901.0195
$ws.Range("J32").Value = 15311.125
$ws.Range("K32").Value = 3901.0195
$ws.Range("L32").Value = 15311.125
$ws.Range("M32").Value = -3614.0195
$ws.Range("N32").Value = -15885.125
$ws.Range("H122").Value = 2937.25
$ws.Range("I122").Value = 999.5
$ws.Range("K122").Value = 2998.5
$ws.Range("M122").Value = -548.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7549.05
$ws.Range("I134").Value = 8239.944
$ws.Range("K134").Value = 24719.832
$ws.Range("M134").Value = -22184.832

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1037.3529
$ws.Range("I22").Value = 426.22223
$ws.Range("J22").Value = 1724.875
$ws.Range("K22").Value = 426.22223
$ws.Range("L22").Value = 1724.875
$ws.Range("M22").Value = -76.22223000000002
$ws.Range("N22").Value = -2424.875
$ws.Range("H58").Value = 1674822.1
$ws.Range("I58").Value = 2901101.5
$ws.Range("J58").Value = 2622.7273
$ws.Range("K58").Value = 2901101.5
$ws.Range("L58").Value = 2622.7273
$ws.Range("M58").Value = -2900898.5
$ws.Range("N58").Value = -3028.7273
$ws.Range("H99").Value = 2575.4
$ws.Range("J99").Value = 2728.5
$ws.Range("L99").Value = 2728.5
$ws.Range("N99").Value = -5724.5
$ws.Range("H122").Value = 4668.6665
$ws.Range("I122").Value = 3500
$ws.Range("J122").Value = 5253
$ws.Range("K122").Value = 10500
$ws.Range("L122").Value = 15759
$ws.Range("M122").Value = -8050
$ws.Range("N122").Value = -20659
$ws.Range("H126").Value = 2575.4
$ws.Range("J126").Value = 2728.5
$ws.Range("L126").Value = 8185.5
$ws.Range("N126").Value = -13125.5
$ws.Range("H136").Value = 1674822.1
$ws.Range("I136").Value = 2901101.5
$ws.Range("J136").Value = 2622.7273
$ws.Range("K136").Value = 8703304.5
$ws.Range("L136").Value = 7868.1819
$ws.Range("M136").Value = -8700754.5
$ws.Range("N136").Value = -12968.1819
$ws.Range("H141").Value = 57673.5
$ws.Range("J141").Value = 55626.855
$ws.Range("L141").Value = 55626.855
$ws.Range("N141").Value = -65986.85500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 624.6111
$ws.Range("I5").Value = 535.9231
$ws.Range("K5").Value = 1607.7693
$ws.Range("M5").Value = -1495.7693
$ws.Range("H44").Value = 1558.3334
$ws.Range("J44").Value = 2200
$ws.Range("L44").Value = 6600
$ws.Range("N44").Value = -7396
$ws.Range("H92").Value = 400
$ws.Range("I92").Value = 400
$ws.Range("K92").Value = 1200
$ws.Range("H122").Value = 834.3333
$ws.Range("I122").Value = 635
$ws.Range("K122").Value = 5715
$ws.Range("M122").Value = -3265
$ws.Range("H135").Value = 624.6111
$ws.Range("I135").Value = 535.9231
$ws.Range("K135").Value = 4823.3079
$ws.Range("M135").Value = -2288.3079
$ws.Range("H139").Value = 4577.6577
$ws.Range("I139").Value = 6426.2383
$ws.Range("J139").Value = 2294.1177
$ws.Range("K139").Value = 19278.7149
$ws.Range("L139").Value = 6882.353099999999
$ws.Range("M139").Value = -14138.7149
$ws.Range("H140").Value = 2362.7273
$ws.Range("J140").Value = 5466
$ws.Range("L140").Value = 16398
$ws.Range("N140").Value = -26758
$ws.Range("H141").Value = 3275.45
$ws.Range("I141").Value = 2673.1177
$ws.Range("K141").Value = 8019.353099999999
$ws.Range("M141").Value = -2839.353099999999
$ws.Range("M92").Value = 48
$ws.Range("N139").Value = -17162.3531

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3927.5715
$ws.Range("I102").Value = 4098.8
$ws.Range("K102").Value = 4098.8
$ws.Range("M102").Value = -2476.8
$ws.Range("H113").Value = 691.8461
$ws.Range("I113").Value = 514.17645
$ws.Range("K113").Value = 514.17645
$ws.Range("M113").Value = 1655.82355

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3128.3333
$ws.Range("I7").Value = 3205.1428
$ws.Range("K7").Value = 3205.1428
$ws.Range("M7").Value = -3093.1428
$ws.Range("H40").Value = 11601.125
$ws.Range("I40").Value = 12184
$ws.Range("K40").Value = 12184
$ws.Range("M40").Value = -12048
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H126").Value = 3128.3333
$ws.Range("I126").Value = 3205.1428
$ws.Range("K126").Value = 9615.428400000001
$ws.Range("M126").Value = -7145.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 43078.75
$ws.Range("J46").Value = 43078.75
$ws.Range("L46").Value = 43078.75
$ws.Range("N46").Value = -43540.75
$ws.Range("H107").Value = 634.2222
$ws.Range("I107").Value = 468.05884
$ws.Range("J107").Value = 916.7
$ws.Range("K107").Value = 1404.17652
$ws.Range("L107").Value = 2750.1
$ws.Range("M107").Value = 515.82348
$ws.Range("N107").Value = -6590.1
$ws.Range("H132").Value = 1169.5853
$ws.Range("I132").Value = 804.32355
$ws.Range("J132").Value = 2943.7144
$ws.Range("K132").Value = 2412.97065
$ws.Range("L132").Value = 8831.143199999999
$ws.Range("M132").Value = 117.0293500000002
$ws.Range("N132").Value = -13891.1432
$ws.Range("H134").Value = 43078.75
$ws.Range("J134").Value = 43078.75
$ws.Range("L134").Value = 129236.25
$ws.Range("N134").Value = -134306.25
